# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the commit: the scraper previously only pulled team statistics;
# this adds the team's season record (W/L/T) as three new trailing columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1): AD1="Wins", AE1="Losses", AF1="Ties" ---
# Reuse the existing header formatting (bold font, thin border, centered)
# from the last pre-existing header cell (AC1) so the new headers match
# the rest of the header row exactly, instead of creating a brand new style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-56): same record repeated for every player on roster ---
$lastRow = 56
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 76   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 86   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
